$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Totals / counters ---
$ws.Range("E11").Value2 = 360533
$ws.Range("F13").Value2 = 7

# --- Insert a new data row before the current last row (row 21) so the
#     table grows from 6 to 7 period rows. Copy row 20 (a "normal" styled
#     row) down into the newly inserted row 21 so it picks up the same
#     formatting as the other interior rows; the old row 21 (with the
#     "last row" styling) shifts down to row 22 automatically. ---
$ws.Rows(21).Insert()
$ws.Range("B20:J20").Copy($ws.Range("B21:J21"))

# Re-number the periods (newest period 2507 on top, oldest 2501 at the
# bottom) and update the corresponding "Valor Mora" amounts.
$ws.Range("E16").Value2 = "2507"
$ws.Range("F16").Value2 = 48533

$ws.Range("E17").Value2 = "2506"
$ws.Range("F17").Value2 = 52000

$ws.Range("E18").Value2 = "2505"
$ws.Range("F18").Value2 = 52000

$ws.Range("E19").Value2 = "2504"
$ws.Range("F19").Value2 = 52000

$ws.Range("E20").Value2 = "2503"
$ws.Range("F20").Value2 = 52000

$ws.Range("E21").Value2 = "2502"
$ws.Range("F21").Value2 = 52000

$ws.Range("E22").Value2 = "2501"
$ws.Range("F22").Value2 = 52000

# NOTE: the single row-insert above already pushed every row below the
# table down by one, so the signature block ("___" underline previously
# on row 26, captions previously on row 27) has already shifted down to
# rows 27/28 - no further action is required there.

Write-Output "done"
